$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Range("C15").Value = 2
$ws.Range("G14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 12
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = -7.692307692307
$ws.Range("L15").Value = 20
$ws.Range("M15").Value = 9.090909090909
$ws.Range("N15").Value = -69.230769230769

# Row 16
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = -77.777777777777
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -27.777777777777
$ws.Range("I16").Value = 127
$ws.Range("J16").Value = 140
$ws.Range("K16").Value = -9.285714285714
$ws.Range("L16").Value = 12.389380530973
$ws.Range("M16").Value = -3.053435114503
$ws.Range("N16").Value = -74.751491053677

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = -45.454545454545
$ws.Range("F17").Value = 30
$ws.Range("G17").Value = 37
$ws.Range("H17").Value = -18.918918918918
$ws.Range("I17").Value = 251
$ws.Range("J17").Value = 276
$ws.Range("K17").Value = -9.057971014492
$ws.Range("L17").Value = 8.658008658008
$ws.Range("M17").Value = 116.379310344828
$ws.Range("N17").Value = -29.494382022471

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 11.111111111111
$ws.Range("I18").Value = 72
$ws.Range("J18").Value = 98
$ws.Range("K18").Value = -26.530612244898
$ws.Range("L18").Value = 7.462686567164
$ws.Range("M18").Value = 67.441860465116
$ws.Range("N18").Value = -71.653543307086

# Row 19
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = 55.172413793103
$ws.Range("I19").Value = 265
$ws.Range("J19").Value = 227
$ws.Range("K19").Value = 16.740088105726
$ws.Range("L19").Value = 9.504132231404
$ws.Range("M19").Value = 77.852348993288
$ws.Range("N19").Value = -25.352112676056

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 48
$ws.Range("J20").Value = 44
$ws.Range("K20").Value = 9.090909090909
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 92
$ws.Range("N20").Value = -79.310344827586

# Row 21
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").Value = 109
$ws.Range("G21").Value = 103
$ws.Range("H21").Value = 5.825242718446
$ws.Range("I21").Value = 777
$ws.Range("J21").Value = 805
$ws.Range("K21").Value = -3.478260869565
$ws.Range("L21").Value = 8.36820083682
$ws.Range("M21").Value = 63.235294117647
$ws.Range("N21").Value = -55.827174530983

# Row 22
$ws.Range("M22").Value = -40

# Row 23
$ws.Range("C23").Value = 11
$ws.Range("D23").Value = 10
$ws.Range("E23").Value = 10
$ws.Range("F23").Value = 42
$ws.Range("G23").Value = 29
$ws.Range("H23").Value = 44.827586206896
$ws.Range("I23").Value = 251
$ws.Range("J23").Value = 240
$ws.Range("K23").Value = 4.583333333333
$ws.Range("L23").Value = 11.555555555555
$ws.Range("M23").Value = 79.285714285714

# Row 24
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = -29.411764705882
$ws.Range("F24").Value = 61
$ws.Range("G24").Value = 54
$ws.Range("H24").Value = 12.962962962963
$ws.Range("I24").Value = 592
$ws.Range("J24").Value = 453
$ws.Range("K24").Value = 30.684326710816
$ws.Range("L24").Value = 8.823529411764
$ws.Range("M24").Value = 57.446808510638

# Row 25
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = 7
$ws.Range("H25").Value = 14.285714285714
$ws.Range("I25").Value = 160
$ws.Range("J25").Value = 92
$ws.Range("K25").Value = 73.91304347826
$ws.Range("L25").Value = -1.234567901234

# Row 26
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 17
$ws.Range("E26").Value = -58.823529411764
$ws.Range("F26").Value = 44
$ws.Range("G26").Value = 55
$ws.Range("H26").Value = -20
$ws.Range("I26").Value = 333
$ws.Range("J26").Value = 460
$ws.Range("K26").Value = -27.608695652173
$ws.Range("L26").Value = -5.93220338983
$ws.Range("M26").Value = -19.565217391304

# Row 27
$ws.Range("C27").Value = 2
$ws.Range("G14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = -25
$ws.Range("I27").Value = 13
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = -27.777777777777
$ws.Range("L27").Value = -31.578947368421

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 2
$ws.Range("G14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = -50
$ws.Range("H14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 25
$ws.Range("I28").Value = 32
$ws.Range("J28").Value = 36
$ws.Range("K28").Value = -11.111111111111
$ws.Range("L28").Value = -8.571428571428

# Row 29
$ws.Range("D29").Value = 1
$ws.Range("G14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = -100
$ws.Range("H14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("G29").Value = 1
$ws.Range("G14").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("H29").Value = -100
$ws.Range("H14").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("J29").Value = 7
$ws.Range("K29").Value = -42.857142857142
$ws.Range("M29").Value = -82.608695652173

# Row 30
$ws.Range("D30").Value = 1
$ws.Range("G14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100
$ws.Range("H14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("G30").Value = 1
$ws.Range("G14").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("H30").Value = -100
$ws.Range("H14").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("J30").Value = 5
$ws.Range("K30").Value = -20
$ws.Range("M30").Value = -80

# Row 31
$ws.Range("F31").NumberFormat = "@"
$ws.Range("F31").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("H31").PasteSpecial(-4122)

# Row 33
$ws.Range("D33").Value = 1
$ws.Range("G14").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("E33").Value = -100
$ws.Range("H14").Copy()
$ws.Range("E33").PasteSpecial(-4122)
$ws.Range("G33").Value = 1
$ws.Range("G14").Copy()
$ws.Range("G33").PasteSpecial(-4122)
$ws.Range("H33").Value = -100
$ws.Range("H14").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("J33").Value = 2

$excel.CutCopyMode = $false